# Auto-generated Excel COM-interop script
$wb = $excel.ActiveWorkbook

# --- 1) Append row 12 (2025-11-12 @ DET) to the four per-player box-score sheets ---
$ws = $wb.Worksheets.Item("Points")
$ws.Cells.Item(12, 1).Value = "'2025-11-12"
$ws.Cells.Item(12, 2).Value = "DET"
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 6
$ws.Cells.Item(12, 7).Value = 12
$ws.Cells.Item(12, 8).Value = 20
$ws.Cells.Item(12, 9).Value = 21
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 9
$ws.Cells.Item(12, 12).Value = 11
$ws.Cells.Item(12, 13).Value = 15
$ws.Cells.Item(12, 14).Value = 7

$ws = $wb.Worksheets.Item("Assists")
$ws.Cells.Item(12, 1).Value = "'2025-11-12"
$ws.Cells.Item(12, 2).Value = "DET"
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 3
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 3
$ws.Cells.Item(12, 7).Value = 2
$ws.Cells.Item(12, 8).Value = 2
$ws.Cells.Item(12, 9).Value = 2
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 5
$ws.Cells.Item(12, 13).Value = 2
$ws.Cells.Item(12, 14).Value = 1

$ws = $wb.Worksheets.Item("Rebounds")
$ws.Cells.Item(12, 1).Value = "'2025-11-12"
$ws.Cells.Item(12, 2).Value = "DET"
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 2
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 6
$ws.Cells.Item(12, 7).Value = 3
$ws.Cells.Item(12, 8).Value = 2
$ws.Cells.Item(12, 9).Value = 14
$ws.Cells.Item(12, 10).Value = 1
$ws.Cells.Item(12, 11).Value = 4
$ws.Cells.Item(12, 12).Value = 2
$ws.Cells.Item(12, 13).Value = 4
$ws.Cells.Item(12, 14).Value = 5

$ws = $wb.Worksheets.Item("3PM")
$ws.Cells.Item(12, 1).Value = "'2025-11-12"
$ws.Cells.Item(12, 2).Value = "DET"
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 3
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 2
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 2
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 13).Value = 2
$ws.Cells.Item(12, 14).Value = 0

# --- 2) Refresh the four "Avg *" leaderboard sheets (new 12-game averages, re-sorted) ---
$ws = $wb.Worksheets.Item("Avg Points")
$ws.Cells.Item(2, 1).Value = "Josh Giddey"
$ws.Cells.Item(2, 2).Value = 21.44444444444444
$ws.Cells.Item(3, 1).Value = "Nikola Vučević"
$ws.Cells.Item(3, 2).Value = 16.09090909090909
$ws.Cells.Item(4, 1).Value = "Matas Buzelis"
$ws.Cells.Item(4, 2).Value = 14.90909090909091
$ws.Cells.Item(5, 1).Value = "Ayo Dosunmu"
$ws.Cells.Item(5, 2).Value = 14.88888888888889
$ws.Cells.Item(6, 1).Value = "Kevin Huerter"
$ws.Cells.Item(6, 2).Value = 14.63636363636364
$ws.Cells.Item(7, 1).Value = "Tre Jones"
$ws.Cells.Item(7, 2).Value = 13.27272727272727
$ws.Cells.Item(8, 1).Value = "Jalen Smith"
$ws.Cells.Item(8, 2).Value = 9.909090909090908
$ws.Cells.Item(9, 1).Value = "Isaac Okoro"
$ws.Cells.Item(9, 2).Value = 8.545454545454545
$ws.Cells.Item(10, 1).Value = "Patrick Williams"
$ws.Cells.Item(10, 2).Value = 8.090909090909092
$ws.Cells.Item(11, 1).Value = "Jevon Carter"
$ws.Cells.Item(11, 2).Value = 5.666666666666667
$ws.Cells.Item(12, 1).Value = "Dalen Terry"
$ws.Cells.Item(12, 2).Value = 1.714285714285714
$ws.Cells.Item(13, 1).Value = "Julian Phillips"
$ws.Cells.Item(13, 2).Value = 0.9

$ws = $wb.Worksheets.Item("Avg Assists")
$ws.Cells.Item(2, 1).Value = "Josh Giddey"
$ws.Cells.Item(2, 2).Value = 9.333333333333334
$ws.Cells.Item(3, 1).Value = "Tre Jones"
$ws.Cells.Item(3, 2).Value = 5.454545454545454
$ws.Cells.Item(4, 1).Value = "Nikola Vučević"
$ws.Cells.Item(4, 2).Value = 3.636363636363636
$ws.Cells.Item(5, 1).Value = "Ayo Dosunmu"
$ws.Cells.Item(5, 2).Value = 3
$ws.Cells.Item(6, 1).Value = "Kevin Huerter"
$ws.Cells.Item(6, 2).Value = 3
$ws.Cells.Item(7, 1).Value = "Jevon Carter"
$ws.Cells.Item(7, 2).Value = 2
$ws.Cells.Item(8, 1).Value = "Isaac Okoro"
$ws.Cells.Item(8, 2).Value = 1.818181818181818
$ws.Cells.Item(9, 1).Value = "Patrick Williams"
$ws.Cells.Item(9, 2).Value = 1.545454545454545
$ws.Cells.Item(10, 1).Value = "Matas Buzelis"
$ws.Cells.Item(10, 2).Value = 1.181818181818182
$ws.Cells.Item(11, 1).Value = "Jalen Smith"
$ws.Cells.Item(11, 2).Value = 1.181818181818182
$ws.Cells.Item(12, 1).Value = "Dalen Terry"
$ws.Cells.Item(12, 2).Value = 0.5714285714285714
$ws.Cells.Item(13, 1).Value = "Julian Phillips"
$ws.Cells.Item(13, 2).Value = 0.1

$ws = $wb.Worksheets.Item("Avg Rebounds")
$ws.Cells.Item(2, 1).Value = "Nikola Vučević"
$ws.Cells.Item(2, 2).Value = 9.909090909090908
$ws.Cells.Item(3, 1).Value = "Josh Giddey"
$ws.Cells.Item(3, 2).Value = 9.555555555555555
$ws.Cells.Item(4, 1).Value = "Jalen Smith"
$ws.Cells.Item(4, 2).Value = 6
$ws.Cells.Item(5, 1).Value = "Matas Buzelis"
$ws.Cells.Item(5, 2).Value = 5.727272727272728
$ws.Cells.Item(6, 1).Value = "Tre Jones"
$ws.Cells.Item(6, 2).Value = 4
$ws.Cells.Item(7, 1).Value = "Kevin Huerter"
$ws.Cells.Item(7, 2).Value = 3.909090909090909
$ws.Cells.Item(8, 1).Value = "Patrick Williams"
$ws.Cells.Item(8, 2).Value = 3.090909090909091
$ws.Cells.Item(9, 1).Value = "Isaac Okoro"
$ws.Cells.Item(9, 2).Value = 2.727272727272727
$ws.Cells.Item(10, 1).Value = "Ayo Dosunmu"
$ws.Cells.Item(10, 2).Value = 2.555555555555555
$ws.Cells.Item(11, 1).Value = "Jevon Carter"
$ws.Cells.Item(11, 2).Value = 1
$ws.Cells.Item(12, 1).Value = "Julian Phillips"
$ws.Cells.Item(12, 2).Value = 0.7
$ws.Cells.Item(13, 1).Value = "Dalen Terry"
$ws.Cells.Item(13, 2).Value = 0.4285714285714285

$ws = $wb.Worksheets.Item("Avg 3PM")
$ws.Cells.Item(2, 1).Value = "Ayo Dosunmu"
$ws.Cells.Item(2, 2).Value = 1.888888888888889
$ws.Cells.Item(3, 1).Value = "Nikola Vučević"
$ws.Cells.Item(3, 2).Value = 1.818181818181818
$ws.Cells.Item(4, 1).Value = "Matas Buzelis"
$ws.Cells.Item(4, 2).Value = 1.818181818181818
$ws.Cells.Item(5, 1).Value = "Josh Giddey"
$ws.Cells.Item(5, 2).Value = 1.666666666666667
$ws.Cells.Item(6, 1).Value = "Jevon Carter"
$ws.Cells.Item(6, 2).Value = 1.666666666666667
$ws.Cells.Item(7, 1).Value = "Kevin Huerter"
$ws.Cells.Item(7, 2).Value = 1.545454545454545
$ws.Cells.Item(8, 1).Value = "Jalen Smith"
$ws.Cells.Item(8, 2).Value = 1.545454545454545
$ws.Cells.Item(9, 1).Value = "Patrick Williams"
$ws.Cells.Item(9, 2).Value = 1.545454545454545
$ws.Cells.Item(10, 1).Value = "Isaac Okoro"
$ws.Cells.Item(10, 2).Value = 1.181818181818182
$ws.Cells.Item(11, 1).Value = "Tre Jones"
$ws.Cells.Item(11, 2).Value = 0.5454545454545454
$ws.Cells.Item(12, 1).Value = "Dalen Terry"
$ws.Cells.Item(12, 2).Value = 0.1428571428571428
$ws.Cells.Item(13, 1).Value = "Julian Phillips"
$ws.Cells.Item(13, 2).Value = 0.1

# --- 3) Add the new "Team Points" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Team Points"

$newSheet.Cells.Item(1, 1).Value = "Game Time (PST)"
$newSheet.Cells.Item(1, 2).Value = "Opponent"
$newSheet.Cells.Item(1, 3).Value = "Team Points"
$newSheet.Cells.Item(1, 4).Value = "Opponent Points"
$newSheet.Cells.Item(1, 5).Value = "Game Total Points"

$headerRange = $newSheet.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$newSheet.Cells.Item(2, 1).Value = "'2025-10-22"
$newSheet.Cells.Item(2, 2).Value = "DET"
$newSheet.Cells.Item(2, 3).Value = 115
$newSheet.Cells.Item(2, 4).Value = 111
$newSheet.Cells.Item(2, 5).Value = 226
$newSheet.Cells.Item(3, 1).Value = "'2025-10-25"
$newSheet.Cells.Item(3, 2).Value = "ORL"
$newSheet.Cells.Item(3, 3).Value = 110
$newSheet.Cells.Item(3, 4).Value = 98
$newSheet.Cells.Item(3, 5).Value = 208
$newSheet.Cells.Item(4, 1).Value = "'2025-10-27"
$newSheet.Cells.Item(4, 2).Value = "ATL"
$newSheet.Cells.Item(4, 3).Value = 128
$newSheet.Cells.Item(4, 4).Value = 123
$newSheet.Cells.Item(4, 5).Value = 251
$newSheet.Cells.Item(5, 1).Value = "'2025-10-29"
$newSheet.Cells.Item(5, 2).Value = "SAC"
$newSheet.Cells.Item(5, 3).Value = 126
$newSheet.Cells.Item(5, 4).Value = 113
$newSheet.Cells.Item(5, 5).Value = 239
$newSheet.Cells.Item(6, 1).Value = "'2025-10-31"
$newSheet.Cells.Item(6, 2).Value = "NYK"
$newSheet.Cells.Item(6, 3).Value = 135
$newSheet.Cells.Item(6, 4).Value = 125
$newSheet.Cells.Item(6, 5).Value = 260
$newSheet.Cells.Item(7, 1).Value = "'2025-11-02"
$newSheet.Cells.Item(7, 2).Value = "NYK"
$newSheet.Cells.Item(7, 3).Value = 116
$newSheet.Cells.Item(7, 4).Value = 128
$newSheet.Cells.Item(7, 5).Value = 244
$newSheet.Cells.Item(8, 1).Value = "'2025-11-04"
$newSheet.Cells.Item(8, 2).Value = "PHI"
$newSheet.Cells.Item(8, 3).Value = 113
$newSheet.Cells.Item(8, 4).Value = 111
$newSheet.Cells.Item(8, 5).Value = 224
$newSheet.Cells.Item(9, 1).Value = "'2025-11-07"
$newSheet.Cells.Item(9, 2).Value = "MIL"
$newSheet.Cells.Item(9, 3).Value = 110
$newSheet.Cells.Item(9, 4).Value = 126
$newSheet.Cells.Item(9, 5).Value = 236
$newSheet.Cells.Item(10, 1).Value = "'2025-11-08"
$newSheet.Cells.Item(10, 2).Value = "CLE"
$newSheet.Cells.Item(10, 3).Value = 122
$newSheet.Cells.Item(10, 4).Value = 128
$newSheet.Cells.Item(10, 5).Value = 250
$newSheet.Cells.Item(11, 1).Value = "'2025-11-10"
$newSheet.Cells.Item(11, 2).Value = "SAS"
$newSheet.Cells.Item(11, 3).Value = 117
$newSheet.Cells.Item(11, 4).Value = 121
$newSheet.Cells.Item(11, 5).Value = 238
$newSheet.Cells.Item(12, 1).Value = "'2025-11-12"
$newSheet.Cells.Item(12, 2).Value = "DET"
$newSheet.Cells.Item(12, 3).Value = 113
$newSheet.Cells.Item(12, 4).Value = 124
$newSheet.Cells.Item(12, 5).Value = 237

Write-Host "Edit complete."
